# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# D/E columns are stored as plain text in the sheet, so numeric-looking Price
# values (single decimal point, e.g. "232.03") are written with a leading
# apostrophe - Excel's standard "store as text" quote-prefix - to stop them
# from being auto-converted to numbers. Values that already contain more
# than one '.' (e.g. "37.311.87") can never parse as a number, so they are
# left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.311.87'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '2.060.72'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''232.03'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('D8').Value = '''57.05'
$ws.Range('E8').Value = '  +0.98%  '
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('D10').Value = '''58.08'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = '''0.0760'
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').Value = '2.364.16'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '''14.60'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '''20.70'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '''0.777'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D17').Value = '''5.13'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').Value = '2.060.84'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').Value = '37.183.77'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').Value = '''6.35'
$ws.Range('E20').Value = '  +7.85%  '
$ws.Range('D21').Value = '''69.31'
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').Value = '''226.20'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('D27').Value = '''166.15'
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('E28').Value = '  +6.39%  '
$ws.Range('D29').Value = '''8.76'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').Value = '''19.05'
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('E32').Value = '  -0.18%  '
$ws.Range('D33').Value = '''4.46'
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').Value = '''4.59'
$ws.Range('E35').Value = '  +5.77%  '
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  -2.04%  '
$ws.Range('D40').Value = '''5.68'
$ws.Range('E40').Value = '  -4.44%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = '1.469.74'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').Value = '''96.16'
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('D44').Value = '''0.0936'
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('D47').Value = '''4.18'
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('D49').Value = '''15.10'
$ws.Range('E49').Value = '  -5.32%  '
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('E51').Value = '  +1.14%  '
